$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at 126:127, shifting existing rows 126-215 down to 128-217
$ws.Rows("126:127").Insert()

# Populate new row 126 (newest week entry, Primera quality)
$ws.Range("A126").Value = 9
$ws.Range("B126").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C126").Value = "Metropolitana"
$ws.Range("D126").Value = 44589
$ws.Range("E126").Value = 13
$ws.Range("F126").Value = 100112017
$ws.Range("G126").Value = "Apio"
$ws.Range("H126").Value = "Americana (o)"
$ws.Range("I126").Value = "Primera"
$ws.Range("J126").Value = 79
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 7000
$ws.Range("M126").Value = 6494
$ws.Range("N126").Value = "$/docena de matas"
$ws.Range("O126").Value = "Región de Coquimbo"
$ws.Range("P126").Value = 1082
$ws.Range("Q126").Value = 6
$ws.Range("R126").Value = "Hortaliza"

# Populate new row 127 (newest week entry, Segunda quality)
$ws.Range("A127").Value = 9
$ws.Range("B127").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C127").Value = "Metropolitana"
$ws.Range("D127").Value = 44589
$ws.Range("E127").Value = 13
$ws.Range("F127").Value = 100112017
$ws.Range("G127").Value = "Apio"
$ws.Range("H127").Value = "Americana (o)"
$ws.Range("I127").Value = "Segunda"
$ws.Range("J127").Value = 25
$ws.Range("K127").Value = 5000
$ws.Range("L127").Value = 5000
$ws.Range("M127").Value = 5000
$ws.Range("N127").Value = "$/docena de matas"
$ws.Range("O127").Value = "Región de Coquimbo"
$ws.Range("P127").Value = 833
$ws.Range("Q127").Value = 6
$ws.Range("R127").Value = "Hortaliza"
